$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 14, 15 rotate (cyclic shift of match data) ---
$ws.Range("B14").Value = 6800432
$ws.Range("E14").Value = "Pontypridd Town"
$ws.Range("F14").Value = "Penybont"
$ws.Range("G14").Value = 0
$ws.Range("I14").Value = "D"
$ws.Range("J14").Value = 3.75
$ws.Range("K14").Value = 3.4
$ws.Range("L14").Value = 1.8
$ws.Range("M14").Value = 4.2
$ws.Range("N14").Value = 3.5
$ws.Range("O14").Value = 1.85
$ws.Range("P14").Value = 0.5
$ws.Range("Q14").Value = 1.825
$ws.Range("R14").Value = 1.975
$ws.Range("T14").Value = 1.975
$ws.Range("U14").Value = 1.825
$ws.Range("V14").Value = -1
$ws.Range("W14").Value = 2.5
$ws.Range("Y14").Value = 0.825
$ws.Range("AB14").Value = 0.825
$ws.Range("B15").Value = 6800431
$ws.Range("E15").Value = "Bala Town"
$ws.Range("F15").Value = "Connahs Quay"
$ws.Range("G15").Value = 1
$ws.Range("I15").Value = "H"
$ws.Range("J15").Value = 2.6
$ws.Range("K15").Value = 3.2
$ws.Range("L15").Value = 2.4
$ws.Range("M15").Value = 2.6
$ws.Range("N15").Value = 3.2
$ws.Range("O15").Value = 2.4
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 1.95
$ws.Range("R15").Value = 1.85
$ws.Range("T15").Value = 2
$ws.Range("U15").Value = 1.8
$ws.Range("V15").Value = 1.6
$ws.Range("W15").Value = -1
$ws.Range("Y15").Value = 0.95
$ws.Range("AB15").Value = 0.8

# --- rows 69, 70, 71 rotate (cyclic shift of match data) ---
$ws.Range("B69").Value = 6800036
$ws.Range("E69").Value = "Caernarfon Town"
$ws.Range("F69").Value = "Aberystwyth"
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 1.4
$ws.Range("K69").Value = 4.5
$ws.Range("L69").Value = 5.75
$ws.Range("M69").Value = 1.363
$ws.Range("N69").Value = 4.75
$ws.Range("O69").Value = 6
$ws.Range("P69").Value = -1.25
$ws.Range("S69").Value = 3
$ws.Range("T69").Value = 1.775
$ws.Range("U69").Value = 2.025
$ws.Range("V69").Value = 0.363
$ws.Range("AA69").Value = 0
$ws.Range("AB69").Value = 0
$ws.Range("B70").Value = 6800472
$ws.Range("E70").Value = "Barry Town"
$ws.Range("F70").Value = "Pontypridd Town"
$ws.Range("G70").Value = 2
$ws.Range("J70").Value = 2.6
$ws.Range("K70").Value = 3.2
$ws.Range("L70").Value = 2.5
$ws.Range("M70").Value = 2.6
$ws.Range("N70").Value = 3.2
$ws.Range("O70").Value = 2.5
$ws.Range("P70").Value = 0
$ws.Range("Q70").Value = 1.975
$ws.Range("R70").Value = 1.825
$ws.Range("S70").Value = 2.5
$ws.Range("T70").Value = 1.95
$ws.Range("U70").Value = 1.85
$ws.Range("V70").Value = 1.6
$ws.Range("Y70").Value = 0.9750000000000001
$ws.Range("AA70").Value = -1
$ws.Range("AB70").Value = 0.8500000000000001
$ws.Range("B71").Value = 6800473
$ws.Range("E71").Value = "Cardiff MU"
$ws.Range("F71").Value = "Newtown"
$ws.Range("H71").Value = 1
$ws.Range("J71").Value = 2.875
$ws.Range("K71").Value = 3.4
$ws.Range("L71").Value = 2.15
$ws.Range("M71").Value = 2.875
$ws.Range("N71").Value = 3.5
$ws.Range("O71").Value = 2.1
$ws.Range("P71").Value = 0.25
$ws.Range("Q71").Value = 1.875
$ws.Range("R71").Value = 1.925
$ws.Range("T71").Value = 1.825
$ws.Range("U71").Value = 1.975
$ws.Range("V71").Value = 1.875
$ws.Range("Y71").Value = 0.875
$ws.Range("AA71").Value = 0.825
$ws.Range("AB71").Value = -1

# --- rows 77, 78 rotate (cyclic shift of match data) ---
$ws.Range("B77").Value = 6800480
$ws.Range("E77").Value = "Penybont"
$ws.Range("F77").Value = "Aberystwyth"
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 2
$ws.Range("I77").Value = "A"
$ws.Range("J77").Value = 1.333
$ws.Range("K77").Value = 4.8
$ws.Range("L77").Value = 6.5
$ws.Range("M77").Value = 1.4
$ws.Range("N77").Value = 4.333
$ws.Range("O77").Value = 5.75
$ws.Range("P77").Value = -1.25
$ws.Range("Q77").Value = 1.925
$ws.Range("R77").Value = 1.875
$ws.Range("S77").Value = 3
$ws.Range("T77").Value = 1.95
$ws.Range("U77").Value = 1.85
$ws.Range("V77").Value = -1
$ws.Range("X77").Value = 4.75
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = 0.875
$ws.Range("AA77").Value = -1
$ws.Range("AB77").Value = 0.8500000000000001
$ws.Range("B78").Value = 6800039
$ws.Range("E78").Value = "Connahs Quay"
$ws.Range("F78").Value = "Caernarfon Town"
$ws.Range("G78").Value = 6
$ws.Range("H78").Value = 1
$ws.Range("I78").Value = "H"
$ws.Range("J78").Value = 1.444
$ws.Range("K78").Value = 4.5
$ws.Range("L78").Value = 5
$ws.Range("M78").Value = 1.333
$ws.Range("N78").Value = 5
$ws.Range("O78").Value = 6.5
$ws.Range("P78").Value = -1.5
$ws.Range("Q78").Value = 1.875
$ws.Range("R78").Value = 1.925
$ws.Range("S78").Value = 3.25
$ws.Range("T78").Value = 1.85
$ws.Range("U78").Value = 1.95
$ws.Range("V78").Value = 0.333
$ws.Range("X78").Value = -1
$ws.Range("Y78").Value = 0.875
$ws.Range("Z78").Value = -1
$ws.Range("AA78").Value = 0.8500000000000001
$ws.Range("AB78").Value = -1

# --- rows 128, 129, 130 rotate (cyclic shift of match data) ---
$ws.Range("B128").Value = 6800518
$ws.Range("E128").Value = "Colwyn Bay"
$ws.Range("F128").Value = "Pontypridd Town"
$ws.Range("H128").Value = 3
$ws.Range("I128").Value = "A"
$ws.Range("J128").Value = 2.2
$ws.Range("K128").Value = 3.75
$ws.Range("L128").Value = 2.7
$ws.Range("M128").Value = 2.25
$ws.Range("N128").Value = 3.8
$ws.Range("O128").Value = 2.6
$ws.Range("P128").Value = -0.25
$ws.Range("T128").Value = 1.825
$ws.Range("U128").Value = 1.975
$ws.Range("W128").Value = -1
$ws.Range("X128").Value = 1.6
$ws.Range("Y128").Value = -1
$ws.Range("Z128").Value = 0.7749999999999999
$ws.Range("AA128").Value = 0.825
$ws.Range("B129").Value = 6800519
$ws.Range("E129").Value = "Haverfordwest County"
$ws.Range("F129").Value = "Bala Town"
$ws.Range("J129").Value = 2.6
$ws.Range("K129").Value = 3.5
$ws.Range("L129").Value = 2.375
$ws.Range("M129").Value = 2.55
$ws.Range("N129").Value = 3.4
$ws.Range("O129").Value = 2.5
$ws.Range("P129").Value = 0
$ws.Range("Q129").Value = 1.9
$ws.Range("R129").Value = 1.9
$ws.Range("S129").Value = 2.5
$ws.Range("T129").Value = 1.975
$ws.Range("U129").Value = 1.825
$ws.Range("X129").Value = 1.5
$ws.Range("Z129").Value = 0.8999999999999999
$ws.Range("AA129").Value = 0.9750000000000001
$ws.Range("B130").Value = 6800053
$ws.Range("E130").Value = "Cardiff MU"
$ws.Range("F130").Value = "Caernarfon Town"
$ws.Range("H130").Value = 2
$ws.Range("I130").Value = "D"
$ws.Range("J130").Value = 2.3
$ws.Range("K130").Value = 3.4
$ws.Range("L130").Value = 2.75
$ws.Range("M130").Value = 2.75
$ws.Range("O130").Value = 2.4
$ws.Range("Q130").Value = 2.025
$ws.Range("R130").Value = 1.775
$ws.Range("S130").Value = 2.75
$ws.Range("T130").Value = 1.9
$ws.Range("U130").Value = 1.9
$ws.Range("W130").Value = 2.4
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 0
$ws.Range("Z130").Value = 0
$ws.Range("AA130").Value = 0.8999999999999999

# --- rows 175, 176 rotate (cyclic shift of match data) ---
$ws.Range("B175").Value = 7721596
$ws.Range("E175").Value = "Bala Town"
$ws.Range("F175").Value = "Newtown"
$ws.Range("G175").Value = 1
$ws.Range("I175").Value = "D"
$ws.Range("J175").Value = 2
$ws.Range("K175").Value = 3.5
$ws.Range("L175").Value = 3
$ws.Range("M175").Value = 2
$ws.Range("O175").Value = 3
$ws.Range("P175").Value = -0.25
$ws.Range("Q175").Value = 1.85
$ws.Range("R175").Value = 1.95
$ws.Range("T175").Value = 1.925
$ws.Range("U175").Value = 1.875
$ws.Range("V175").Value = -1
$ws.Range("W175").Value = 2.4
$ws.Range("Y175").Value = -0.5
$ws.Range("Z175").Value = 0.475
$ws.Range("AA175").Value = -1
$ws.Range("AB175").Value = 0.875
$ws.Range("B176").Value = 7721623
$ws.Range("E176").Value = "Haverfordwest County"
$ws.Range("F176").Value = "Barry Town"
$ws.Range("G176").Value = 2
$ws.Range("I176").Value = "H"
$ws.Range("J176").Value = 1.95
$ws.Range("K176").Value = 3.25
$ws.Range("L176").Value = 3.4
$ws.Range("M176").Value = 1.666
$ws.Range("O176").Value = 4.333
$ws.Range("P176").Value = -0.75
$ws.Range("Q176").Value = 1.925
$ws.Range("R176").Value = 1.875
$ws.Range("T176").Value = 1.975
$ws.Range("U176").Value = 1.825
$ws.Range("V176").Value = 0.6659999999999999
$ws.Range("W176").Value = -1
$ws.Range("Y176").Value = 0.4625
$ws.Range("Z176").Value = -0.5
$ws.Range("AA176").Value = 0.9750000000000001
$ws.Range("AB176").Value = -1
